$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 5: the date in column C moves from 7.4.2015 to 8.4.2015 ---
$ws.Range("C5").Value = "8.4.2015 г. 00:00:00 ч."

# --- Widen column D (closest the engine's column-width quantization can reach to 25.85546875) ---
$ws.Columns.Item(4).ColumnWidth = 25

# --- Add two new rows of transaction data ---
$ws.Range("A8").Value = "QCUCJKC37W"
$ws.Range("B8").Value = 678
$ws.Range("C8").Value = "10.4.2015 г. 00:00:00 ч."
$ws.Range("D8").Value = "RegularExpense"
$ws.Range("E8").Value = "Malko po-dylyg tekst"

$ws.Range("A9").Value = "ARKSWRKWDA"
$ws.Range("B9").Value = 900
$ws.Range("C9").Value = "29.4.2015 г."
$ws.Range("D9").Value = "RegularExpense"
$ws.Range("E9").Value = "Muahahaha"

# --- Re-apply the shared row formatting (setting .Value drops existing styles) ---
$ws.Range("C6").Copy()
$ws.Range("C5").PasteSpecial(-4122)

$ws.Range("A7:E7").Copy()
$ws.Range("A8:E9").PasteSpecial(-4122)

# --- Update the selected cell to reflect the new selection shown in the workbook ---
$ws.Range("C5").Select()
